$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (row 11) correct-answer count
$ws.Range("B11").Value = 5

# Update Total row (row 12): corrected total and corr/total marks text
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
